# Weekly update for the Alcachofa (artichoke) price sheet.
# - Inserts two new report rows (7 and 8) with this week's market entries.
# - Shifts the previously-existing rows 7-15 down to rows 9-17.
# - Row 17 (previously the old row 15 data) is refreshed with a new date
#   and updated price figures.
# - A new row 18 is appended carrying the original (pre-shift) row 15 data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Insert two blank rows at position 7, pushing old rows 7-15 to 9-17 ---
$ws.Range("A7:A8").EntireRow.Insert()

function Set-DataRow {
    param($Row, $Fecha, $Variedad, $Volumen, $PrecioMin, $PrecioMax, $PrecioProm, $Unidad, $PrecioKg, $KgUnidades)

    $ws.Cells.Item($Row, 1).Value = 11
    $ws.Cells.Item($Row, 2).Value = "Vega Monumental Concepción"
    $ws.Cells.Item($Row, 3).Value = "Bíobío"
    $ws.Cells.Item($Row, 4).Value = $Fecha
    $ws.Cells.Item($Row, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
    $ws.Cells.Item($Row, 5).Value = 8
    $ws.Cells.Item($Row, 6).Value = 100112013
    $ws.Cells.Item($Row, 7).Value = "Alcachofa"
    $ws.Cells.Item($Row, 8).Value = $Variedad
    $ws.Cells.Item($Row, 9).Value = "Primera"
    $ws.Cells.Item($Row, 10).Value = $Volumen
    $ws.Cells.Item($Row, 11).Value = $PrecioMin
    $ws.Cells.Item($Row, 12).Value = $PrecioMax
    $ws.Cells.Item($Row, 13).Value = $PrecioProm
    $ws.Cells.Item($Row, 14).Value = $Unidad
    $ws.Cells.Item($Row, 15).Value = "Provincia de Limarí"
    $ws.Cells.Item($Row, 16).Value = $PrecioKg
    $ws.Cells.Item($Row, 17).Value = $KgUnidades
    $ws.Cells.Item($Row, 18).Value = "Hortaliza"
}

# --- New row 7: Española, $/caja 30 unidades ---
Set-DataRow 7 44421 "Española" 100 14000 15000 14500 "`$/caja 30 unidades" 483 30

# --- New row 8: Madrigal, $/caja 40 unidades ---
Set-DataRow 8 44426 "Madrigal" 50 12000 13000 12600 "`$/caja 40 unidades" 315 40

# --- Row 17 (was old row 15, shifted down by the insert) gets refreshed data ---
Set-DataRow 17 44420 "Española" 100 14000 15000 14500 "`$/caja 30 unidades" 483 30

# --- New row 18 appended, carrying the original row 15 data forward ---
Set-DataRow 18 44376 "Española" 100 19000 20000 19500 "`$/caja 30 unidades" 650 30
